$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9358363151550293
$ws.Range("B1").Value = 2.763909339904785
$ws.Range("C1").Value = 4.214782238006592
$ws.Range("D1").Value = 0.9318075180053711
$ws.Range("E1").Value = 0.8098978400230408
